$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the pipeline / pipeline source values to use underscores instead of
# spaces (parameter-friendly names for the build process).
$ws.Range("A2").Value = "Pipeline_1"
$ws.Range("A3").Value = "Pipeline_2"
$ws.Range("B2").Value = "Pipeline_Source_1"
$ws.Range("B3").Value = "Pipeline_Source_2"

$ws.Range("F7").Select()
